# Fruta / hortaliza, semanal
# Insert 3 new weekly report rows before the current row 485, shifting the
# existing rows 485-493 down to 486-496 (net new dimension A1:R496).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows at 485 (pushes old 485:493 down to 488:496).
$ws.Range("A485:R487").Insert()

# New row 485
$ws.Cells.Item(485, 1).Value = 9
$ws.Cells.Item(485, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(485, 3).Value = "Metropolitana"
$ws.Cells.Item(485, 4).Value = 44628
$ws.Cells.Item(485, 5).Value = 13
$ws.Cells.Item(485, 6).Value = 100112024
$ws.Cells.Item(485, 7).Value = "Choclo"
$ws.Cells.Item(485, 8).Value = "Choclero"
$ws.Cells.Item(485, 9).Value = "Primera"
$ws.Cells.Item(485, 10).Value = 6100
$ws.Cells.Item(485, 11).Value = 160
$ws.Cells.Item(485, 12).Value = 180
$ws.Cells.Item(485, 13).Value = 170
$ws.Cells.Item(485, 14).Value = "$/unidad"
$ws.Cells.Item(485, 15).Value = "Región Metropolitana"
$ws.Cells.Item(485, 16).Value = 170
$ws.Cells.Item(485, 17).Value = 1
$ws.Cells.Item(485, 18).Value = "Hortaliza"

# New row 486
$ws.Cells.Item(486, 1).Value = 9
$ws.Cells.Item(486, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(486, 3).Value = "Metropolitana"
$ws.Cells.Item(486, 4).Value = 44628
$ws.Cells.Item(486, 5).Value = 13
$ws.Cells.Item(486, 6).Value = 100112024
$ws.Cells.Item(486, 7).Value = "Choclo"
$ws.Cells.Item(486, 8).Value = "Choclero"
$ws.Cells.Item(486, 9).Value = "Primera"
$ws.Cells.Item(486, 10).Value = 7900
$ws.Cells.Item(486, 11).Value = 160
$ws.Cells.Item(486, 12).Value = 180
$ws.Cells.Item(486, 13).Value = 170
$ws.Cells.Item(486, 14).Value = "$/unidad"
$ws.Cells.Item(486, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(486, 16).Value = 170
$ws.Cells.Item(486, 17).Value = 1
$ws.Cells.Item(486, 18).Value = "Hortaliza"

# New row 487
$ws.Cells.Item(487, 1).Value = 9
$ws.Cells.Item(487, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(487, 3).Value = "Metropolitana"
$ws.Cells.Item(487, 4).Value = 44628
$ws.Cells.Item(487, 5).Value = 13
$ws.Cells.Item(487, 6).Value = 100112024
$ws.Cells.Item(487, 7).Value = "Choclo"
$ws.Cells.Item(487, 8).Value = "Dulce o Americano"
$ws.Cells.Item(487, 9).Value = "Primera"
$ws.Cells.Item(487, 10).Value = 4300
$ws.Cells.Item(487, 11).Value = 130
$ws.Cells.Item(487, 12).Value = 150
$ws.Cells.Item(487, 13).Value = 140
$ws.Cells.Item(487, 14).Value = "$/unidad"
$ws.Cells.Item(487, 15).Value = "Región Metropolitana"
$ws.Cells.Item(487, 16).Value = 140
$ws.Cells.Item(487, 17).Value = 1
$ws.Cells.Item(487, 18).Value = "Hortaliza"
